# configuratie.xlsx - apply the "conceptversie" intro-text update.
#
# Summary of the change:
#   - intro_tekst!B4 ("Dit tabellenboek is een conceptversie.") gets expanded
#     with extra wording about [naam] / [subsetnaam].
#   - A new row is inserted right after it (new row 5) with a second new
#     sentence about [subregio].
#   - The old "Iedere GGD kan ..." row shifts down (was row 5, now row 6),
#     and the trailing blank "tekst" row shifts down too (was row 6, now row 7).
#   - intro_tekst becomes the active / selected sheet (it picks up
#     tabSelected, replacing "algemeen"), with B5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intro_tekst")

# Insert a fresh row above the current row 5 ("Iedere GGD kan ...") so the
# new sentence about [subregio] gets its own row, and everything below
# shifts down by one.
$ws.Rows("5:5").Insert()

# Expand the existing conceptversie sentence in place (still row 4).
$ws.Range("B4").Value = "Dit tabellenboek is een conceptversie. Je kunt [naam] gebruiken voor de hoofdsubset, of andere subsets aanduiden met [subsetnaam]."

# Fill the newly inserted row 5 with the [subregio] sentence (column A keeps
# using the same "tekst" label as the surrounding rows).
$ws.Range("A5").Value = "tekst"
$ws.Range("B5").Value = "Je zou bijvoorbeeld de subregio aan kunnen geven met [subregio]."

# Make intro_tekst the active sheet/tab with B5 selected, matching the
# workbook's saved view state after the edit.
$ws.Activate()
[void]$ws.Range("B5").Select()
